$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table data (rows 2-11, columns A:C) replacing the old rows 2-15.
$data = @(
    @("Login_Verification",        "All type of login execution",                         "Y"),
    @("Navigate_Verification",     "All Page UI and Navigation verification testing ",     "N"),
    @("Teacher_Verification",      "Teacher Module Testing ",                              "N"),
    @("Parent_Verification",       "Parent Module Testing ",                               "N"),
    @("Student_Verification",      "Student Module Testing ",                              "N"),
    @("User_Verification",         "User Module Testing ",                                 "N"),
    @("Resources_Verification",    "Resources Module Testing ",                            "N"),
    @("SimpleSearch_Verification", "Simple Search  Module Testing ",                       "N"),
    @("AdvanceSearch_Verification","Advance  Search Module Testing ",                      "N"),
    @("Class_Verification",        "Class Verification Module Testing ",                   "N")
)

# Remove the rows that will no longer be needed (old sheet had data through row 15,
# new sheet only needs through row 11).
$ws.Rows("12:15").Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("B18").Select() | Out-Null
